# "Generate Report for Handoff"
# Updates the localization-status report: the "Handed back: in sync with
# en-US" status becomes "Ready for handoff" (now the shortest value in its
# columns, so those columns are narrowed), and the handoff/generate
# timestamps tied to that status refresh are bumped forward a few minutes.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "Ready for handoff"

# --- Timestamps refreshed alongside the status change ---
$ws1.Range("G2").Value = "2016-09-03 21:02:20"
$ws3.Range("H2").Value = "2016-09-03 21:02:20"
$ws2.Range("H2").Value = "2016-09-03 21:02:15"

# --- Column widths shrink to fit the new, shorter status text ---
$ws1.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws1.Columns.Item(6).ColumnWidth = 16.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws3.Columns.Item(3).ColumnWidth = 16.333333333333332
